$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.868.17"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "3.438.21"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.60"
$ws.Range("E5").Value = "  +4.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.42"
$ws.Range("E6").Value = "  +8.07%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "3.434.27"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.643"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.01"
$ws.Range("E12").Value = "  +6.24%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.43"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "3.982.19"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.84"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "3.432.11"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "66.825.84"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.03"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.34"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.42"
$ws.Range("E23").Value = "  +10.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.99"
$ws.Range("E24").Value = "  +20.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.34"
$ws.Range("E25").Value = "  +6.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.94"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.96"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.94"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.98"
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.05"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  +12.18%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "64.54"
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "597.16"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.146"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").Value = "3.192.11"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +5.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0429"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  +5.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.135"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  +20.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +5.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.59"
$ws.Range("E51").Value = "  +3.46%  "
